$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for columns C (runs), D (balls), E (fours), F (sixes)
# for rows 2 through 13, as described by the diff.
$data = @{
    2  = @("17","22","1","0")
    3  = @("41","22","7","1")
    4  = @("30","23","1","2")
    5  = @("19","13","2","1")
    6  = @("13","11","2","0")
    7  = @("32","27","3","1")
    8  = @("9","4","2","0")
    10 = @("5","9","0","0")
    11 = @("6","2","0","1")
    12 = @("2","7","0","0")
    13 = @("4","9","0","0")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
}
